# adicionando exemplos de turmas na planilha matriculas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: nomes das turmas (horario)
$ws.Range("A1").Value = "segqua1112"
$ws.Range("A2").Value = "terqui0910"
$ws.Range("A3").Value = "segqua2122"
$ws.Range("A4").Value = "terqui1314"
$ws.Range("A5").Value = "segqua0708"
$ws.Range("A6").Value = "terqui1415"

# Column B: exemplos de alunos matriculados em cada turma
$ws.Range("B1").Value = "manel, bernardo, pipico"
$ws.Range("B2").Value = "ana, alice, andré"
$ws.Range("B3").Value = "thales, amanda, letícia"
$ws.Range("B4").Value = "natasha, luisa, gabriela"
$ws.Range("B5").Value = "mariana, davi, eduardo"
$ws.Range("B6").Value = "rafael, jennifer, luiz"

# Marcadores auxiliares (células sublinhadas, sem valor) usados na planilha
$ws.Range("M13").Font.Underline = $true
$ws.Range("C16").Font.Underline = $true

# Seleção final conforme o estado salvo do arquivo
$ws.Range("M13").Select()
